$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy full-row formatting from the last existing data row (row 89) down to the new rows (90-97)
$ws.Range("A89:V89").Copy()
$ws.Range("A90:A97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 90 (Indice 89)
$ws.Cells.Item(90,1).Value2 = 89
$ws.Cells.Item(90,2).Value2 = "belgium"
$ws.Cells.Item(90,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(90,4).Value2 = "2023-2024"
$ws.Cells.Item(90,5).Value2 = 45226.86458333334
$ws.Cells.Item(90,6).Value2 = "St. Truiden"
$ws.Cells.Item(90,7).Value2 = 2
$ws.Cells.Item(90,8).Value2 = "RWDM"
$ws.Cells.Item(90,9).Value2 = 1
$ws.Cells.Item(90,10).Value2 = 1.98
$ws.Cells.Item(90,11).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(90,12).Value2 = 1.83
$ws.Cells.Item(90,13).Value2 = "27/10/2023 20:13"
$ws.Cells.Item(90,14).Value2 = 3.69
$ws.Cells.Item(90,15).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(90,16).Value2 = 3.75
$ws.Cells.Item(90,17).Value2 = "27/10/2023 20:38"
$ws.Cells.Item(90,18).Value2 = 3.74
$ws.Cells.Item(90,19).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(90,20).Value2 = 4.36
$ws.Cells.Item(90,21).Value2 = "27/10/2023 19:49"
$ws.Cells.Item(90,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/st-truiden-rwd-molenbeek/4zUWKOdi/"

# Row 91 (Indice 90)
$ws.Cells.Item(91,1).Value2 = 90
$ws.Cells.Item(91,2).Value2 = "belgium"
$ws.Cells.Item(91,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(91,4).Value2 = "2023-2024"
$ws.Cells.Item(91,5).Value2 = 45227.66666666666
$ws.Cells.Item(91,6).Value2 = "KV Mechelen"
$ws.Cells.Item(91,7).Value2 = 0
$ws.Cells.Item(91,8).Value2 = "Cercle Brugge KSV"
$ws.Cells.Item(91,9).Value2 = 2
$ws.Cells.Item(91,10).Value2 = 3.07
$ws.Cells.Item(91,11).Value2 = "22/10/2023 16:12"
$ws.Cells.Item(91,12).Value2 = 3.72
$ws.Cells.Item(91,13).Value2 = "28/10/2023 15:30"
$ws.Cells.Item(91,14).Value2 = 3.69
$ws.Cells.Item(91,15).Value2 = "22/10/2023 16:12"
$ws.Cells.Item(91,16).Value2 = 4.05
$ws.Cells.Item(91,17).Value2 = "28/10/2023 15:30"
$ws.Cells.Item(91,18).Value2 = 2.15
$ws.Cells.Item(91,19).Value2 = "22/10/2023 16:12"
$ws.Cells.Item(91,20).Value2 = 1.93
$ws.Cells.Item(91,21).Value2 = "28/10/2023 15:20"
$ws.Cells.Item(91,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/kv-mechelen-cercle-brugge/hd0E2ssN/"

# Row 92 (Indice 91)
$ws.Cells.Item(92,1).Value2 = 91
$ws.Cells.Item(92,2).Value2 = "belgium"
$ws.Cells.Item(92,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(92,4).Value2 = "2023-2024"
$ws.Cells.Item(92,5).Value2 = 45227.76041666666
$ws.Cells.Item(92,6).Value2 = "Eupen"
$ws.Cells.Item(92,7).Value2 = 2
$ws.Cells.Item(92,8).Value2 = "Charleroi"
$ws.Cells.Item(92,9).Value2 = 0
$ws.Cells.Item(92,10).Value2 = 3.12
$ws.Cells.Item(92,11).Value2 = "21/10/2023 21:12"
$ws.Cells.Item(92,12).Value2 = 3.88
$ws.Cells.Item(92,13).Value2 = "28/10/2023 18:13"
$ws.Cells.Item(92,14).Value2 = 3.71
$ws.Cells.Item(92,15).Value2 = "21/10/2023 21:12"
$ws.Cells.Item(92,16).Value2 = 4.03
$ws.Cells.Item(92,17).Value2 = "28/10/2023 18:13"
$ws.Cells.Item(92,18).Value2 = 2.2
$ws.Cells.Item(92,19).Value2 = "21/10/2023 21:12"
$ws.Cells.Item(92,20).Value2 = 1.88
$ws.Cells.Item(92,21).Value2 = "28/10/2023 18:12"
$ws.Cells.Item(92,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/eupen-charleroi/OMeA3NRG/"

# Row 93 (Indice 92)
$ws.Cells.Item(93,1).Value2 = 92
$ws.Cells.Item(93,2).Value2 = "belgium"
$ws.Cells.Item(93,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(93,4).Value2 = "2023-2024"
$ws.Cells.Item(93,5).Value2 = 45227.86458333334
$ws.Cells.Item(93,6).Value2 = "Anderlecht"
$ws.Cells.Item(93,7).Value2 = 5
$ws.Cells.Item(93,8).Value2 = "Leuven"
$ws.Cells.Item(93,9).Value2 = 1
$ws.Cells.Item(93,10).Value2 = 1.58
$ws.Cells.Item(93,11).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(93,12).Value2 = 1.45
$ws.Cells.Item(93,13).Value2 = "28/10/2023 20:40"
$ws.Cells.Item(93,14).Value2 = 4.36
$ws.Cells.Item(93,15).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(93,16).Value2 = 4.84
$ws.Cells.Item(93,17).Value2 = "28/10/2023 20:44"
$ws.Cells.Item(93,18).Value2 = 4.81
$ws.Cells.Item(93,19).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(93,20).Value2 = 7.08
$ws.Cells.Item(93,21).Value2 = "28/10/2023 20:44"
$ws.Cells.Item(93,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/anderlecht-leuven/pWd643CA/"

# Row 94 (Indice 93)
$ws.Cells.Item(94,1).Value2 = 93
$ws.Cells.Item(94,2).Value2 = "belgium"
$ws.Cells.Item(94,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(94,4).Value2 = "2023-2024"
$ws.Cells.Item(94,5).Value2 = 45228.5625
$ws.Cells.Item(94,6).Value2 = "Club Brugge KV"
$ws.Cells.Item(94,7).Value2 = 2
$ws.Cells.Item(94,8).Value2 = "Antwerp"
$ws.Cells.Item(94,9).Value2 = 1
$ws.Cells.Item(94,10).Value2 = 1.93
$ws.Cells.Item(94,11).Value2 = "22/10/2023 14:42"
$ws.Cells.Item(94,12).Value2 = 1.88
$ws.Cells.Item(94,13).Value2 = "29/10/2023 13:28"
$ws.Cells.Item(94,14).Value2 = 3.88
$ws.Cells.Item(94,15).Value2 = "22/10/2023 14:42"
$ws.Cells.Item(94,16).Value2 = 3.75
$ws.Cells.Item(94,17).Value2 = "29/10/2023 13:28"
$ws.Cells.Item(94,18).Value2 = 3.71
$ws.Cells.Item(94,19).Value2 = "22/10/2023 14:42"
$ws.Cells.Item(94,20).Value2 = 4.25
$ws.Cells.Item(94,21).Value2 = "29/10/2023 13:23"
$ws.Cells.Item(94,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/club-brugge-antwerp/27715qc4/"

# Row 95 (Indice 94)
$ws.Cells.Item(95,1).Value2 = 94
$ws.Cells.Item(95,2).Value2 = "belgium"
$ws.Cells.Item(95,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(95,4).Value2 = "2023-2024"
$ws.Cells.Item(95,5).Value2 = 45228.66666666666
$ws.Cells.Item(95,6).Value2 = "Westerlo"
$ws.Cells.Item(95,7).Value2 = 1
$ws.Cells.Item(95,8).Value2 = "Royale Union SG"
$ws.Cells.Item(95,9).Value2 = 3
$ws.Cells.Item(95,10).Value2 = 4.12
$ws.Cells.Item(95,11).Value2 = "22/10/2023 17:12"
$ws.Cells.Item(95,12).Value2 = 5.92
$ws.Cells.Item(95,13).Value2 = "29/10/2023 15:52"
$ws.Cells.Item(95,14).Value2 = 4.07
$ws.Cells.Item(95,15).Value2 = "22/10/2023 17:12"
$ws.Cells.Item(95,16).Value2 = 4.36
$ws.Cells.Item(95,17).Value2 = "29/10/2023 15:52"
$ws.Cells.Item(95,18).Value2 = 1.74
$ws.Cells.Item(95,19).Value2 = "22/10/2023 17:12"
$ws.Cells.Item(95,20).Value2 = 1.56
$ws.Cells.Item(95,21).Value2 = "29/10/2023 15:52"
$ws.Cells.Item(95,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/westerlo-royale-union-sg/xUyocL4p/"

# Row 96 (Indice 95)
$ws.Cells.Item(96,1).Value2 = 95
$ws.Cells.Item(96,2).Value2 = "belgium"
$ws.Cells.Item(96,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(96,4).Value2 = "2023-2024"
$ws.Cells.Item(96,5).Value2 = 45228.77083333334
$ws.Cells.Item(96,6).Value2 = "Gent"
$ws.Cells.Item(96,7).Value2 = 3
$ws.Cells.Item(96,8).Value2 = "St. Liege"
$ws.Cells.Item(96,9).Value2 = 1
$ws.Cells.Item(96,10).Value2 = 1.48
$ws.Cells.Item(96,11).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(96,12).Value2 = 1.75
$ws.Cells.Item(96,13).Value2 = "29/10/2023 18:25"
$ws.Cells.Item(96,14).Value2 = 4.9
$ws.Cells.Item(96,15).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(96,16).Value2 = 4.14
$ws.Cells.Item(96,17).Value2 = "29/10/2023 18:25"
$ws.Cells.Item(96,18).Value2 = 5.94
$ws.Cells.Item(96,19).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(96,20).Value2 = 4.47
$ws.Cells.Item(96,21).Value2 = "29/10/2023 18:25"
$ws.Cells.Item(96,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/gent-st-liege/jXukduKj/"

# Row 97 (Indice 96)
$ws.Cells.Item(97,1).Value2 = 96
$ws.Cells.Item(97,2).Value2 = "belgium"
$ws.Cells.Item(97,3).Value2 = "jupiler-pro-league"
$ws.Cells.Item(97,4).Value2 = "2023-2024"
$ws.Cells.Item(97,5).Value2 = 45228.80208333334
$ws.Cells.Item(97,6).Value2 = "Kortrijk"
$ws.Cells.Item(97,7).Value2 = 0
$ws.Cells.Item(97,8).Value2 = "Genk"
$ws.Cells.Item(97,9).Value2 = 3
$ws.Cells.Item(97,10).Value2 = 4.5
$ws.Cells.Item(97,11).Value2 = "22/10/2023 20:42"
$ws.Cells.Item(97,12).Value2 = 6.07
$ws.Cells.Item(97,13).Value2 = "29/10/2023 18:51"
$ws.Cells.Item(97,14).Value2 = 4.48
$ws.Cells.Item(97,15).Value2 = "22/10/2023 20:42"
$ws.Cells.Item(97,16).Value2 = 4.77
$ws.Cells.Item(97,17).Value2 = "29/10/2023 18:51"
$ws.Cells.Item(97,18).Value2 = 1.67
$ws.Cells.Item(97,19).Value2 = "22/10/2023 20:42"
$ws.Cells.Item(97,20).Value2 = 1.51
$ws.Cells.Item(97,21).Value2 = "29/10/2023 18:51"
$ws.Cells.Item(97,22).Value2 = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/kortrijk-genk/WnaI11dT/"

